$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.241.54"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "3.680.07"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "684.33"
$ws.Range("E5").Value = "  -3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.51"
$ws.Range("E6").Value = "  -5.61%  "
$ws.Range("D7").Value = "3.680.02"
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("E10").Value = "  -8.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.39"
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("E13").Value = "  -5.82%  "
$ws.Range("E14").Value = "  -6.80%  "
$ws.Range("D15").Value = "4.301.74"
$ws.Range("E15").Value = "  -3.59%  "
$ws.Range("D16").Value = "3.679.46"
$ws.Range("E16").Value = "  -3.52%  "
$ws.Range("D17").Value = "69.329.25"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  -6.67%  "
$ws.Range("E20").Value = "  -7.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.87"
$ws.Range("E21").Value = "  -6.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  -7.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.667"
$ws.Range("E23").Value = "  -8.13%  "
$ws.Range("E24").Value = "  -5.22%  "
$ws.Range("D25").Value = "3.827.13"
$ws.Range("E25").Value = "  -3.58%  "
$ws.Range("E26").Value = "  -10.19%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.42"
$ws.Range("E28").Value = "  -5.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  -8.97%  "
$ws.Range("E30").Value = "  -11.06%  "
$ws.Range("E31").Value = "  -11.22%  "
$ws.Range("E32").Value = "  -7.92%  "
$ws.Range("E33").Value = "  -8.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.15"
$ws.Range("E34").Value = "  -7.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.166"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "3.650.67"
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.48"
$ws.Range("E38").Value = "  -7.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.36"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0932"
$ws.Range("E41").Value = "  -8.42%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.950"
$ws.Range("E44").Value = "  -7.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.27"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.38"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").Value = "  -14.31%  "
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  -9.23%  "
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("E51").Value = "  -3.20%  "
